# Edit script for digicode.xlsx
# - Removes the two "ChatGPT Pro" rows (Compartido / Privado) since the
#   shop no longer sells that tier, shifting all following rows up by 2.
# - Updates the surviving "ChatGPT PLUS" rows (Compartido / Privado) with
#   their new subscription detail, price (soles) and product image.
# - Fills in the previously empty MARCA (brand/category) column for the
#   tail of the list (PicsArt ... Windows 10 LTSC).
# - Converts the CODIGO column (A) from static text to a live formula
#   that derives the code from the row number: DIG100, DIG101, ...
# - Restores the view (scroll position / selection) to match where the
#   author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the two discontinued "ChatGPT Pro" rows (rows 18 and 19).
#    Everything below shifts up by two rows automatically.
# ---------------------------------------------------------------------
$ws.Rows("18:19").Delete()

# ---------------------------------------------------------------------
# 2. Refresh the two remaining "ChatGPT PLUS" rows with their new
#    detail / price / image.
# ---------------------------------------------------------------------
# Row 16: ChatGPT PLUS - Compartido
$ws.Range("G16").Value = "Sucripcion x 30 dias"
$ws.Range("H16").Value = 19
$ws.Range("I16").Value = "https://nguyenpremium.com/wp-content/uploads/2024/09/ChatGPT-Plus-1.jpg"

# Row 17: ChatGPT PLUS - Privado
$ws.Range("G17").Value = "Sucripcion x 30 dias"
$ws.Range("H17").Value = 80
$ws.Range("I17").Value = "https://nguyenpremium.com/wp-content/uploads/2024/09/ChatGPT-Plus-1.jpg"

# ---------------------------------------------------------------------
# 3. Fill in the MARCA (column E) values that were missing for the tail
#    of the table (rows 42-52 after the deletion above).
# ---------------------------------------------------------------------
$ws.Range("E42").Value = "DISEÑO"            # PicsArt
$ws.Range("E43").Value = "PRODUCTIVIDAD"     # Power ISO
$ws.Range("E44").Value = "STREAMING"         # PRIME VIDEO
$ws.Range("E45").Value = "INGENIERIA"        # SolidWork 2024
$ws.Range("E46").Value = "INGENIERIA"        # SolidWork 2025
$ws.Range("E47").Value = "PRODUCTIVIDAD"     # UDEMY
$ws.Range("E48").Value = "PRODUCTIVIDAD"     # vTubeGo
$ws.Range("E49").Value = "STREAMING"         # Youtube Premium
$ws.Range("E50").Value = "SISTEMA"           # Microsoft Windows 10 Pro
$ws.Range("E51").Value = "SISTEMA"           # Microsoft Windows 11 Pro
$ws.Range("E52").Value = "SISTEMA"           # Microsoft Windows 10 LTSC

# ---------------------------------------------------------------------
# 4. Turn the CODIGO column into a calculated formula column.
#    A2 gets its own formula, A3:A52 share one formula group (this is
#    how Excel naturally splits it when filled as two separate writes).
# ---------------------------------------------------------------------
$ws.Range("A2").Formula = '="DIG"&TEXT(ROW()-1+100,"000")'
$ws.Range("A3:A52").Formula = '="DIG"&TEXT(ROW()-1+100,"000")'

# Keep the table's calculated-column formula definition in sync too.
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.ListColumns.Item("CODIGO").Range.Formula = '="DIG"&TEXT(ROW()-1+100,"000")'

# ---------------------------------------------------------------------
# 5. Restore the view: scrolled down so row 35 is at the top, with
#    E50:E52 selected (the newly completed SISTEMA brand cells).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E50:E52").Select()
